$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 84
$ws.Range("F3").Value = 138
$ws.Range("F4").Value = 1763
$ws.Range("F5").Value = 3319
$ws.Range("F6").Value = 1021
$ws.Range("F7").Value = 2187
$ws.Range("F8").Value = 2098
$ws.Range("F10").Value = 600
$ws.Range("F12").Value = 1665
$ws.Range("F13").Value = 391
$ws.Range("F15").Value = 41
$ws.Range("F17").Value = 200
$ws.Range("F18").Value = 1571
$ws.Range("F19").Value = 622
$ws.Range("F20").Value = 709
$ws.Range("F21").Value = 595
$ws.Range("F22").Value = 12195
$ws.Range("F23").Value = 12228
$ws.Range("F24").Value = 907
$ws.Range("F27").Value = 28
$ws.Range("F28").Value = 17
$ws.Range("F30").Value = 1915
$ws.Range("F32").Value = 571

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 119
$ws.Range("F7").Value = 19

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 84
$ws.Range("F4").Value = 138
$ws.Range("F5").Value = 1763
$ws.Range("F6").Value = 3319
$ws.Range("F7").Value = 1021
$ws.Range("F8").Value = 2187
$ws.Range("F9").Value = 2098
$ws.Range("F11").Value = 600
$ws.Range("F13").Value = 1665
$ws.Range("F14").Value = 391
$ws.Range("F17").Value = 41
$ws.Range("F21").Value = 200
$ws.Range("F22").Value = 1571
$ws.Range("F23").Value = 622
$ws.Range("F24").Value = 709
$ws.Range("F25").Value = 595
$ws.Range("F26").Value = 12195
$ws.Range("F27").Value = 12228
$ws.Range("F28").Value = 907
$ws.Range("F31").Value = 28
$ws.Range("F32").Value = 17
$ws.Range("F34").Value = 1915
$ws.Range("F35").Value = 119
$ws.Range("F38").Value = 571
$ws.Range("F39").Value = 19
